# Applies the "Started work on SizeStructure1" edit to the W22 schedule sheet.
#
# What happened conceptually: a new "Size Structure" sub-topic was split into
# two parts ("Size Structure I" / "Size Structure II"), which pushed the
# remaining topics (Weight-Length, Condition, Bag Limits, Size Limits,
# Quotas, Coarse Woody Habitat) down by one row (rows 20-26 -> 21-27).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("W22")

# Row 19 becomes the first "Size Structure" session.
$ws.Range("D19").Value = "Size Structure I"
$ws.Range("E19").Value = "SizeStructure1"

# A new row is effectively inserted: the old contents of rows 20-26
# (D:E) shift down to rows 21-27, and row 20 gets the new second
# "Size Structure" session.
$ws.Range("D20").Value = "Size Structure II"
$ws.Range("E20").Value = "SizeStructure2"

$ws.Range("D21").Value = "Weight-Length"
$ws.Range("E21").Value = "WeightLength"

$ws.Range("D22").Value = "Condition"
$ws.Range("E22").Value = "Condition"

$ws.Range("D23").Value = "Bag Limits"
$ws.Range("E23").Value = "BagLimits"

$ws.Range("D24").Value = "Size Limits"
$ws.Range("E24").Value = "SizeLimits"

$ws.Range("D25").Value = "Quotas"
$ws.Range("E25").Value = "Quotas"

$ws.Range("D27").Value = "Coarse Woody Habitat"
$ws.Range("E27").Value = "CWH"

# Update the selection left behind by the editor.
$ws.Range("D30").Select()
